{"js": "// \"Tercera version de cambios\":\n// Append a new paragraph, in Courier New, after the existing trailing\n// (empty) paragraph at the end of the document body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document currently ends with an empty paragraph right before the\n// section break; insert the new paragraph right after it (i.e. at the\n// very end of the body).\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst newParagraph = lastParagraph.insertParagraph(\n  \"Realizando una peque\u00f1a modificaci\u00f3n en cervantes.docx.\",\n  Word.InsertLocation.after\n);\n\n// Give the new text the same Courier New font used elsewhere in the\n// document.\nnewParagraph.font.name = \"Courier New\";\n\nawait context.sync();\n", "ps1": "# \"Tercera version de cambios\":\n# Append a new paragraph, in Courier New, after the existing trailing\n# (empty) paragraph at the end of the document.\n\n$d = $word.ActiveDocument\n\n# Collapse to the very end of the document and add a brand-new paragraph\n# there (after the existing empty trailing paragraph).\n$endRange = $d.Range($d.Content.End, $d.Content.End)\n$endRange.InsertParagraphAfter()\n\n# The paragraph we just created is now the last paragraph in the document.\n$newPara = $d.Paragraphs.Last\n\n# Apply Courier New to the (still empty) paragraph first, so the paragraph\n# mark itself - and whatever gets typed next - inherits the font.\n$newPara.Range.Font.Name = \"Courier New\"\n$newPara.Range.Font.NameBi = \"Courier New\"\n\n# Type the sentence into the new paragraph.\n$insertPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)\n$insertPoint.Text = \"Realizando una peque\u00f1a modificaci\u00f3n en cervantes.docx.\"\n\n# Re-assert Courier New (Latin + complex-script) on the whole paragraph so\n# the run that now holds the text is explicitly formatted, matching the\n# rest of the document's formatting style.\n$newPara.Range.Font.Name = \"Courier New\"\n$newPara.Range.Font.NameBi = \"Courier New\"\n"}
